$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 8: "Improved Non Random Access (iNRA)"
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Add($p.Slides.Count + 1, 2)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "Improved Non Random Access (iNRA)"

$tr8 = $s8.Shapes.Item(2).TextFrame.TextRange
$tr8.Text = "Improvements to the NRA algorithm"
[void]$tr8.InsertAfter("`rUse Length Boundedness to prune search space. ")
[void]$tr8.InsertAfter("`r")

# ---------------------------------------------------------------------------
# Slide 9: "Shortest-First Algorithm"
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Add($p.Slides.Count + 1, 2)
$s9.Shapes.Item(1).TextFrame.TextRange.Text = "Shortest-First Algorithm"

$tr9 = $s9.Shapes.Item(2).TextFrame.TextRange
$tr9.Text = "Scan list in decreasing idf order"
[void]$tr9.InsertAfter("`rTokens with lower idfs examined first")
[void]$tr9.InsertAfter("`rOccur in fewer strings")
[void]$tr9.InsertAfter("`rSmaller lists of strings that they occur in")
[void]$tr9.InsertAfter("`rSmallest chance of false positive candidates")
[void]$tr9.InsertAfter("`r")
[void]$tr9.InsertAfter("`rFind shortest list of candidate strings. ")

$tr9.Paragraphs(2).IndentLevel = 2
$tr9.Paragraphs(3).IndentLevel = 3
$tr9.Paragraphs(4).IndentLevel = 3
$tr9.Paragraphs(5).IndentLevel = 3
$tr9.Paragraphs(6).IndentLevel = 3

# ---------------------------------------------------------------------------
# Slide 10: "Hybrid Algorithm"
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Add($p.Slides.Count + 1, 2)
$s10.Shapes.Item(1).TextFrame.TextRange.Text = "Hybrid Algorithm"
$s10.Shapes.Item(2).TextFrame.TextRange.Text = "ADD STUFF HERE!"

# ---------------------------------------------------------------------------
# Slide 11: "Conclusions"
# ---------------------------------------------------------------------------
$s11 = $p.Slides.Add($p.Slides.Count + 1, 2)
$s11.Shapes.Item(1).TextFrame.TextRange.Text = "Conclusions"

$tr11 = $s11.Shapes.Item(2).TextFrame.TextRange
$tr11.Text = "In General, SF had the best performance of the mentioned algorithms"
[void]$tr11.InsertAfter("`rOnly in certain edge cases did the guarantees of iNRA improve performance")
[void]$tr11.InsertAfter("`rAll algorithms benefit from Length Bounded pruning")
[void]$tr11.InsertAfter("`r")

$tr11.Paragraphs(2).IndentLevel = 2
